$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $newText, 2) | Out-Null
}

function Find-Paragraph($searchText, $wholeWord) {
    $f = $d.Content.Find
    $f.ClearFormatting()
    $f.Text = $searchText
    $f.Forward = $true
    $f.Wrap = 1
    $f.MatchWholeWord = $wholeWord
    $f.MatchWildcards = $false
    $f.Execute() | Out-Null
    $pt = $d.Range($f.Parent.Start, $f.Parent.Start)
    return $pt.Paragraphs(1)
}

function Insert-ParagraphAfterText($anchorText, $newParaText) {
    # Locate the (now single-run) paragraph whose whole text is $anchorText,
    # split a fresh paragraph mark after it, then re-locate that same
    # anchor again (stale object refs from before the split can misreport
    # Start/End, especially for the last paragraph in the body) before
    # grabbing .Next() and filling in its text.
    $anchor = Find-Paragraph $anchorText $true
    $anchor.Range.InsertParagraphAfter()
    $anchor = Find-Paragraph $anchorText $true
    $newPara = $anchor.Next()
    $newPara.Range.Text = $newParaText
}

# --- "Verktyg:" section ---
# "Grafiklösning ännu inte bestämd" -> "Visual Studio som IDE"
# then insert a new paragraph "SFML för grafiken" right after it.
Replace-Text "Grafiklösning ännu inte bestämd" "Visual Studio som IDE"
Insert-ParagraphAfterText "Visual Studio som IDE" "SFML för grafiken"

# --- SPRINT 2 section ---
# "Nära fungerande spel" -> "Grafik"
# then insert a new paragraph "Lägga ut enheter" right after it.
Replace-Text "Nära fungerande spel" "Grafik"
Insert-ParagraphAfterText "Grafik" "Lägga ut enheter"

# --- SPRINT 3 section ---
# "Grafik, UI och möjligtvis AI" -> "Flytta enheter"
# then insert a new paragraph "AI" right after it.
Replace-Text "Grafik, UI och möjligtvis AI" "Flytta enheter"
Insert-ParagraphAfterText "Flytta enheter" "AI"
